# Generate Report for Handoff
# Updates the handoff-generation timestamps and marks the priority ("ht")
# for the batch of files that just got a handoff report generated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date for the newly-handed-off batch (rows 7-12).
# Overview!G7:G12 and de-de!H7:H12 share this same timestamp string.
$overview.Range("G7:G12").Value = "2016-09-06 08:27:38"

# Latest Handoff Datetime on the zh-cn sheet for the same batch.
$zhcn.Range("H7:H12").Value = "2016-09-06 08:27:32"

# Mark Priority as "ht" for the same batch of rows on both locale sheets.
$zhcn.Range("E7:E12").Value = "ht"
$dede.Range("E7:E12").Value = "ht"
